$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(111,8).Value = 2199.8  # H111: 1841.5 -> 2199.8
$ws.Cells.Item(111,9).Value = 3000  # I111: 2016.6666 -> 3000
$ws.Cells.Item(111,11).Value = 9000  # K111: 6049.9998 -> 9000
$ws.Cells.Item(111,13).Value = -5933  # M111: -2982.9998 -> -5933
$ws.Cells.Item(112,8).Value = 1163.0555  # H112: 1161.2858 -> 1163.0555
$ws.Cells.Item(112,10).Value = 1143.2354  # J112: 1144.35 -> 1143.2354
$ws.Cells.Item(112,12).Value = 3429.7062  # L112: 3433.05 -> 3429.7062
$ws.Cells.Item(112,14).Value = -5645.706200000001  # N112: -5649.049999999999 -> -5645.706200000001
$ws.Cells.Item(132,8).Value = 1827.3334  # H132: 1842 -> 1827.3334
$ws.Cells.Item(132,9).Value = 1798.5  # I132: 1799 -> 1798.5
$ws.Cells.Item(132,11).Value = 5395.5  # K132: 5397 -> 5395.5
$ws.Cells.Item(132,13).Value = -2865.5  # M132: -2867 -> -2865.5
$ws.Cells.Item(138,8).Value = 3597.02  # H138: 3658.11 -> 3597.02
$ws.Cells.Item(138,10).Value = 3877.9578  # J138: 3964 -> 3877.9578
$ws.Cells.Item(138,12).Value = 11633.8734  # L138: 11892 -> 11633.8734
$ws.Cells.Item(138,14).Value = -21913.8734  # N138: -22172 -> -21913.8734

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32,8).Value = 3982.3171  # H32: 3757.15 -> 3982.3171
$ws.Cells.Item(32,10).Value = 12994.333  # J32: 12997 -> 12994.333
$ws.Cells.Item(32,12).Value = 12994.333  # L32: 12997 -> 12994.333
$ws.Cells.Item(32,14).Value = -13568.333  # N32: -13571 -> -13568.333
$ws.Cells.Item(45,8).Value = 1736.5714  # H45: 1800.1538 -> 1736.5714
$ws.Cells.Item(45,9).Value = 1717.8462  # I45: 1785.1666 -> 1717.8462
$ws.Cells.Item(45,11).Value = 1717.8462  # K45: 1785.1666 -> 1717.8462
$ws.Cells.Item(45,13).Value = -1340.8462  # M45: -1408.1666 -> -1340.8462
$ws.Cells.Item(88,8).Value = 1877.1111  # H88: 1985.1428 -> 1877.1111
$ws.Cells.Item(88,10).Value = 1877.1111  # J88: 1985.1428 -> 1877.1111
$ws.Cells.Item(88,12).Value = 1877.1111  # L88: 1985.1428 -> 1877.1111
$ws.Cells.Item(88,14).Value = -2689.1111  # N88: -2797.1428 -> -2689.1111
$ws.Cells.Item(91,8).Value = 1877.1111  # H91: 1985.1428 -> 1877.1111
$ws.Cells.Item(91,10).Value = 1877.1111  # J91: 1985.1428 -> 1877.1111
$ws.Cells.Item(91,12).Value = 1877.1111  # L91: 1985.1428 -> 1877.1111
$ws.Cells.Item(91,14).Value = -4685.1111  # N91: -4793.1428 -> -4685.1111
$ws.Cells.Item(122,8).Value = 1593.3334  # H122: 1929.7142 -> 1593.3334
$ws.Cells.Item(122,9).Value = 1048.5714  # I122: 1301.6 -> 1048.5714
$ws.Cells.Item(122,11).Value = 3145.7142  # K122: 3904.8 -> 3145.7142
$ws.Cells.Item(122,13).Value = -695.7142000000003  # M122: -1454.8 -> -695.7142000000003
$ws.Cells.Item(132,8).Value = 2141  # H132: 3505.5 -> 2141
$ws.Cells.Item(132,9).Value = 2141  # I132: 3505.5 -> 2141
$ws.Cells.Item(132,11).Value = 6423  # K132: 10516.5 -> 6423
$ws.Cells.Item(132,13).Value = -3893  # M132: -7986.5 -> -3893

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94,8).Value = 427.54544  # H94: 466.4 -> 427.54544
$ws.Cells.Item(94,9).Value = 430.3  # I94: 473.77777 -> 430.3
$ws.Cells.Item(94,11).Value = 430.3  # K94: 473.77777 -> 430.3
$ws.Cells.Item(94,13).Value = 20.69999999999999  # M94: -22.77776999999998 -> 20.69999999999999
$ws.Cells.Item(105,8).Value = 2874.75  # H105: 2908.875 -> 2874.75
$ws.Cells.Item(105,9).Value = 2836  # I105: 2889.4614 -> 2836
$ws.Cells.Item(105,10).Value = 2991  # J105: 2993 -> 2991
$ws.Cells.Item(105,11).Value = 2836  # K105: 2889.4614 -> 2836
$ws.Cells.Item(105,12).Value = 2991  # L105: 2993 -> 2991
$ws.Cells.Item(105,13).Value = -1089  # M105: -1142.4614 -> -1089
$ws.Cells.Item(105,14).Value = -6485  # N105: -6487 -> -6485
$ws.Cells.Item(107,8).Value = 3218.4546  # H107: 3291.4546 -> 3218.4546
$ws.Cells.Item(107,9).Value = 2739  # I107: 3174.125 -> 2739
$ws.Cells.Item(107,10).Value = 8013  # J107: 3604.3333 -> 8013
$ws.Cells.Item(107,11).Value = 2739  # K107: 3174.125 -> 2739
$ws.Cells.Item(107,12).Value = 8013  # L107: 3604.3333 -> 8013
$ws.Cells.Item(107,13).Value = -819  # M107: -1254.125 -> -819
$ws.Cells.Item(107,14).Value = -11853  # N107: -7444.3333 -> -11853
$ws.Cells.Item(134,8).Value = 1250  # H134: 2678.8333 -> 1250
$ws.Cells.Item(134,9).Value = 1250  # I134: 2678.8333 -> 1250
$ws.Cells.Item(134,11).Value = 3750  # K134: 8036.499899999999 -> 3750
$ws.Cells.Item(134,13).Value = -1215  # M134: -5501.499899999999 -> -1215

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16,8).Value = 1895.5  # H16: 2427.7273 -> 1895.5
$ws.Cells.Item(16,9).Value = 1859  # I16: 2421.5557 -> 1859
$ws.Cells.Item(16,10).Value = 2053.6667  # J16: 2455.5 -> 2053.6667
$ws.Cells.Item(16,11).Value = 1859  # K16: 2421.5557 -> 1859
$ws.Cells.Item(16,12).Value = 2053.6667  # L16: 2455.5 -> 2053.6667
$ws.Cells.Item(16,13).Value = -1572  # M16: -2134.5557 -> -1572
$ws.Cells.Item(16,14).Value = -2627.6667  # N16: -3029.5 -> -2627.6667
$ws.Cells.Item(31,8).Value = 5829.4443  # H31: 6295.625 -> 5829.4443
$ws.Cells.Item(31,9).Value = 4470.5713  # I31: 4899 -> 4470.5713
$ws.Cells.Item(31,10).Value = 6694.1816  # J31: 7133.6 -> 6694.1816
$ws.Cells.Item(31,11).Value = 4470.5713  # K31: 4899 -> 4470.5713
$ws.Cells.Item(31,12).Value = 6694.1816  # L31: 7133.6 -> 6694.1816
$ws.Cells.Item(31,13).Value = -4175.5713  # M31: -4604 -> -4175.5713
$ws.Cells.Item(31,14).Value = -7284.1816  # N31: -7723.6 -> -7284.1816
$ws.Cells.Item(34,8).Value = 5829.4443  # H34: 6295.625 -> 5829.4443
$ws.Cells.Item(34,9).Value = 4470.5713  # I34: 4899 -> 4470.5713
$ws.Cells.Item(34,10).Value = 6694.1816  # J34: 7133.6 -> 6694.1816
$ws.Cells.Item(34,11).Value = 4470.5713  # K34: 4899 -> 4470.5713
$ws.Cells.Item(34,12).Value = 6694.1816  # L34: 7133.6 -> 6694.1816
$ws.Cells.Item(34,13).Value = -4268.5713  # M34: -4697 -> -4268.5713
$ws.Cells.Item(34,14).Value = -7098.1816  # N34: -7537.6 -> -7098.1816
$ws.Cells.Item(113,8).Value = 1895.5  # H113: 2427.7273 -> 1895.5
$ws.Cells.Item(113,9).Value = 1859  # I113: 2421.5557 -> 1859
$ws.Cells.Item(113,10).Value = 2053.6667  # J113: 2455.5 -> 2053.6667
$ws.Cells.Item(113,11).Value = 1859  # K113: 2421.5557 -> 1859
$ws.Cells.Item(113,12).Value = 2053.6667  # L113: 2455.5 -> 2053.6667
$ws.Cells.Item(113,13).Value = 311  # M113: -251.5556999999999 -> 311
$ws.Cells.Item(113,14).Value = -6393.6667  # N113: -6795.5 -> -6393.6667
$ws.Cells.Item(122,8).Value = 1066.8889  # H122: 1190.2 -> 1066.8889
$ws.Cells.Item(122,9).Value = 1143.6666  # I122: 1058.8572 -> 1143.6666
$ws.Cells.Item(122,10).Value = 913.3333  # J122: 1496.6666 -> 913.3333
$ws.Cells.Item(122,11).Value = 3430.9998  # K122: 3176.5716 -> 3430.9998
$ws.Cells.Item(122,12).Value = 2739.9999  # L122: 4489.9998 -> 2739.9999
$ws.Cells.Item(122,13).Value = -980.9998000000001  # M122: -726.5715999999998 -> -980.9998000000001
$ws.Cells.Item(122,14).Value = -7639.9999  # N122: -9389.9998 -> -7639.9999
$ws.Cells.Item(141,8).Value = 321064.47  # H141: 293295.94 -> 321064.47
$ws.Cells.Item(141,10).Value = 321064.47  # J141: 293295.94 -> 321064.47
$ws.Cells.Item(141,12).Value = 321064.47  # L141: 293295.94 -> 321064.47
$ws.Cells.Item(141,14).Value = -331424.47  # N141: -303655.94 -> -331424.47

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23,8).Value = 575.0526  # H23: 602.55554 -> 575.0526
$ws.Cells.Item(23,9).Value = 412.07144  # I23: 437.6154 -> 412.07144
$ws.Cells.Item(23,11).Value = 1236.21432  # K23: 1312.8462 -> 1236.21432
$ws.Cells.Item(23,13).Value = -1001.21432  # M23: -1077.8462 -> -1001.21432
$ws.Cells.Item(32,8).Value = 0  # H32: 50 -> 0
$ws.Cells.Item(32,9).Value = 0  # I32: 50 -> 0
$ws.Cells.Item(32,11).Value = 0  # K32: 150 -> 0
$ws.Cells.Item(32,13).ClearContents()  # M32: was 133
$ws.Cells.Item(34,8).Value = 500  # H34: 578.3333 -> 500
$ws.Cells.Item(34,9).Value = 500  # I34: 578.3333 -> 500
$ws.Cells.Item(34,11).Value = 1500  # K34: 1734.9999 -> 1500
$ws.Cells.Item(34,13).Value = -1416  # M34: -1650.9999 -> -1416
$ws.Cells.Item(39,8).Value = 7924.75  # H39: 8231.666999999999 -> 7924.75
$ws.Cells.Item(39,10).Value = 7924.75  # J39: 8231.666999999999 -> 7924.75
$ws.Cells.Item(39,12).Value = 23774.25  # L39: 24695.001 -> 23774.25
$ws.Cells.Item(39,14).Value = -24362.25  # N39: -25283.001 -> -24362.25
$ws.Cells.Item(55,8).Value = 5331.3335  # H55: 2749.75 -> 5331.3335
$ws.Cells.Item(55,9).Value = 3999.5  # I55: 2833 -> 3999.5
$ws.Cells.Item(55,10).Value = 7995  # J55: 2500 -> 7995
$ws.Cells.Item(55,11).Value = 11998.5  # K55: 8499 -> 11998.5
$ws.Cells.Item(55,12).Value = 23985  # L55: 7500 -> 23985
$ws.Cells.Item(55,13).Value = -11821.5  # M55: -8322 -> -11821.5
$ws.Cells.Item(55,14).Value = -24339  # N55: -7854 -> -24339
$ws.Cells.Item(122,8).Value = 400  # H122: 316.66666 -> 400
$ws.Cells.Item(122,9).Value = 100  # I122: 125 -> 100
$ws.Cells.Item(122,11).Value = 900  # K122: 1125 -> 900
$ws.Cells.Item(122,13).Value = 1550  # M122: 1325 -> 1550
$ws.Cells.Item(136,8).Value = 5833  # H136: 5999.6665 -> 5833
$ws.Cells.Item(136,10).Value = 4500  # J136: 5000 -> 4500
$ws.Cells.Item(136,12).Value = 13500  # L136: 15000 -> 13500
$ws.Cells.Item(136,14).Value = -23700  # N136: -25200 -> -23700

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70,8).Value = 5460.1816  # H70: 5221.5625 -> 5460.1816
$ws.Cells.Item(70,9).Value = 4874  # I70: 4686.25 -> 4874
$ws.Cells.Item(70,10).Value = 5795.143  # J70: 5756.875 -> 5795.143
$ws.Cells.Item(70,11).Value = 4874  # K70: 4686.25 -> 4874
$ws.Cells.Item(70,12).Value = 5795.143  # L70: 5756.875 -> 5795.143
$ws.Cells.Item(70,13).Value = -4604  # M70: -4416.25 -> -4604
$ws.Cells.Item(70,14).Value = -6335.143  # N70: -6296.875 -> -6335.143
$ws.Cells.Item(73,8).Value = 5460.1816  # H73: 5221.5625 -> 5460.1816
$ws.Cells.Item(73,9).Value = 4874  # I73: 4686.25 -> 4874
$ws.Cells.Item(73,10).Value = 5795.143  # J73: 5756.875 -> 5795.143
$ws.Cells.Item(73,11).Value = 4874  # K73: 4686.25 -> 4874
$ws.Cells.Item(73,12).Value = 5795.143  # L73: 5756.875 -> 5795.143
$ws.Cells.Item(73,13).Value = -3938  # M73: -3750.25 -> -3938
$ws.Cells.Item(73,14).Value = -7667.143  # N73: -7628.875 -> -7667.143
$ws.Cells.Item(97,8).Value = 736.5  # H97: 914.9 -> 736.5
$ws.Cells.Item(97,9).Value = 715.5  # I97: 763.2857 -> 715.5
$ws.Cells.Item(97,10).Value = 799.5  # J97: 1268.6666 -> 799.5
$ws.Cells.Item(97,11).Value = 715.5  # K97: 763.2857 -> 715.5
$ws.Cells.Item(97,12).Value = 799.5  # L97: 1268.6666 -> 799.5
$ws.Cells.Item(97,13).Value = -219.5  # M97: -267.2857 -> -219.5
$ws.Cells.Item(97,14).Value = -1791.5  # N97: -2260.6666 -> -1791.5
$ws.Cells.Item(113,8).Value = 3999.4  # H113: 4000 -> 3999.4
$ws.Cells.Item(113,9).Value = 3999  # I113: 4000 -> 3999
$ws.Cells.Item(113,11).Value = 3999  # K113: 4000 -> 3999
$ws.Cells.Item(113,13).Value = -1829  # M113: -1830 -> -1829
$ws.Cells.Item(122,8).Value = 2244.6667  # H122: 1851.1 -> 2244.6667
$ws.Cells.Item(122,9).Value = 1817.5  # I122: 1526.5 -> 1817.5
$ws.Cells.Item(122,10).Value = 3099  # J122: 3149.5 -> 3099
$ws.Cells.Item(122,11).Value = 5452.5  # K122: 4579.5 -> 5452.5
$ws.Cells.Item(122,12).Value = 9297  # L122: 9448.5 -> 9297
$ws.Cells.Item(122,13).Value = -3002.5  # M122: -2129.5 -> -3002.5
$ws.Cells.Item(122,14).Value = -14197  # N122: -14348.5 -> -14197
$ws.Cells.Item(139,8).Value = 70000  # H139: 0 -> 70000
$ws.Cells.Item(139,10).Value = 70000  # J139: 0 -> 70000
$ws.Cells.Item(139,12).Value = 70000  # L139: 0 -> 70000

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16,8).Value = 1201.2858  # H16: 1262 -> 1201.2858
$ws.Cells.Item(16,9).Value = 1000  # I16: 1141.6666 -> 1000
$ws.Cells.Item(16,11).Value = 1000  # K16: 1141.6666 -> 1000
$ws.Cells.Item(16,13).Value = -830  # M16: -971.6666 -> -830
$ws.Cells.Item(61,8).Value = 6993.7085  # H61: 7080.8696 -> 6993.7085
$ws.Cells.Item(61,9).Value = 5792.75  # I61: 5835.0527 -> 5792.75
$ws.Cells.Item(61,11).Value = 5792.75  # K61: 5835.0527 -> 5792.75
$ws.Cells.Item(61,13).Value = -5590.75  # M61: -5633.0527 -> -5590.75
$ws.Cells.Item(82,8).Value = 1535.7273  # H82: 1535.909 -> 1535.7273
$ws.Cells.Item(82,9).Value = 1374.25  # I82: 1374.5 -> 1374.25
$ws.Cells.Item(82,11).Value = 1374.25  # K82: 1374.5 -> 1374.25
$ws.Cells.Item(82,13).Value = -1013.25  # M82: -1013.5 -> -1013.25
$ws.Cells.Item(85,8).Value = 1535.7273  # H85: 1535.909 -> 1535.7273
$ws.Cells.Item(85,9).Value = 1374.25  # I85: 1374.5 -> 1374.25
$ws.Cells.Item(85,11).Value = 1374.25  # K85: 1374.5 -> 1374.25
$ws.Cells.Item(85,13).Value = -126.25  # M85: -126.5 -> -126.25
$ws.Cells.Item(113,8).Value = 6993.7085  # H113: 7080.8696 -> 6993.7085
$ws.Cells.Item(113,9).Value = 5792.75  # I113: 5835.0527 -> 5792.75
$ws.Cells.Item(113,11).Value = 5792.75  # K113: 5835.0527 -> 5792.75
$ws.Cells.Item(113,13).Value = -3622.75  # M113: -3665.0527 -> -3622.75
$ws.Cells.Item(122,8).Value = 6679.9395  # H122: 6835.625 -> 6679.9395
$ws.Cells.Item(122,9).Value = 6476.7915  # I122: 6614.3335 -> 6476.7915
$ws.Cells.Item(122,10).Value = 7221.6665  # J122: 7499.5 -> 7221.6665
$ws.Cells.Item(122,11).Value = 19430.3745  # K122: 19843.0005 -> 19430.3745
$ws.Cells.Item(122,12).Value = 21664.9995  # L122: 22498.5 -> 21664.9995
$ws.Cells.Item(122,13).Value = -16980.3745  # M122: -17393.0005 -> -16980.3745
$ws.Cells.Item(122,14).Value = -26564.9995  # N122: -27398.5 -> -26564.9995
$ws.Cells.Item(136,8).Value = 6620.6665  # H136: 6637.3335 -> 6620.6665
$ws.Cells.Item(136,9).Value = 6736.4287  # I136: 6754.2856 -> 6736.4287
$ws.Cells.Item(136,11).Value = 20209.2861  # K136: 20262.8568 -> 20209.2861
$ws.Cells.Item(136,13).Value = -17659.2861  # M136: -17712.8568 -> -17659.2861

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107,8).Value = 270.88235  # H107: 287.375 -> 270.88235
$ws.Cells.Item(107,9).Value = 207.16667  # I107: 225.36363 -> 207.16667
$ws.Cells.Item(107,11).Value = 621.50001  # K107: 676.0908899999999 -> 621.50001
$ws.Cells.Item(107,13).Value = 1298.49999  # M107: 1243.90911 -> 1298.49999
$ws.Cells.Item(113,8).Value = 235.85715  # H113: 261.16666 -> 235.85715
$ws.Cells.Item(113,9).Value = 108.5  # I113: 113.4 -> 108.5
$ws.Cells.Item(113,11).Value = 325.5  # K113: 340.2 -> 325.5
$ws.Cells.Item(113,13).Value = 1844.5  # M113: 1829.8 -> 1844.5
